$d = $word.ActiveDocument

# The run currently reads:
#   "Biggest by far: when running, also generate a json of the database
#    diffs. Then Electron GUI showing grid"
# It needs to become two runs (same character formatting on both):
#   "Biggest by far: when running, also generate a "
#   "HTMLs of the diffs"

$oldTail = "json of the database diffs. Then Electron GUI showing grid"
$newTail = "HTMLs of the diffs"

# 1) Locate the trailing text that must change and give it its own run
#    by toggling a character-formatting property on and back off again.
#    This lets Word split the run at that boundary without altering any
#    visible formatting (the color stays inherited from the original
#    run's rPr on both halves).
$full = $d.Content.Text
$idx = $full.IndexOf($oldTail)
if ($idx -lt 0) {
    throw "Could not find the target sentence to edit"
}
$startPos = $idx
$endPos = $idx + $oldTail.Length
$splitRange = $d.Range($startPos, $endPos)
$splitRange.Bold = 1
$splitRange.Bold = 0

# 2) Now that the tail is an isolated run, replace its text in place.
$replaceRange = $d.Range($startPos, $endPos)
$found = $replaceRange.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)
if (-not $found) {
    throw "Could not replace the target sentence"
}

# 3) The text replacement can re-merge the edited run back with its
#    identically-formatted neighbour, so force the split once more on
#    the freshly-inserted text to make sure it ends up as its own run.
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf($newTail)
if ($idx2 -lt 0) {
    throw "Could not find the replacement text after the edit"
}
$startPos2 = $idx2
$endPos2 = $idx2 + $newTail.Length
$finalSplitRange = $d.Range($startPos2, $endPos2)
$finalSplitRange.Bold = 1
$finalSplitRange.Bold = 0

Write-Host "Edit applied"
